$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly re-sync reshuffled which calendar day each Higo price record
# belongs to. Rows keep all their descriptive columns (market, product,
# quality, unit, origin, etc.) but the Fecha / Volumen / Precio columns
# (D, M, N, O, P, S) move between rows along two 4-row cycles and two
# simple row swaps:
#   2 -> 13 -> 4 -> 15 -> 2   (row R receives the old D/M/N/O/P/S of the
#                              row that follows it in the cycle)
#   3 -> 14 -> 5 -> 16 -> 3
#   11 <-> 17
#   12 <-> 18

$cycles = @(
    @(2, 13, 4, 15),
    @(3, 14, 5, 16),
    @(11, 17),
    @(12, 18)
)

$cols = @("D", "M", "N", "O", "P", "S")

foreach ($cycle in $cycles) {
    $n = $cycle.Count

    # Snapshot the current (pre-edit) values for every row in this cycle,
    # for every affected column, before writing anything back.
    $snapshot = @{}
    for ($i = 0; $i -lt $n; $i++) {
        $row = $cycle[$i]
        foreach ($col in $cols) {
            $snapshot["$col-$row"] = $ws.Range("$col$row").Value2()
        }
    }

    # Row at position i takes on the old values that belonged to the row
    # at position i+1 (wrapping around), matching the rotation above.
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $n]
        foreach ($col in $cols) {
            $ws.Range("$col$destRow").Value = $snapshot["$col-$srcRow"]
        }
    }
}
